# Update the two time-slot labels in column C (rows 8 and 9) and move the
# active selection on the sheet from C18 to C17, matching the authored
# workbook edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "13:25-13:30"
$ws.Range("C9").Value = "13:30-13:35"

$ws.Range("C17").Select()
